$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: header/field names ---
# H: BombIdentity -> BombNumber
$ws.Range("H1").Value = "BombNumber"
# J: OperatingMode -> Programme
$ws.Range("J1").Value = "Programme"
# L: Energy -> SamplePortion (ml); M: SamplePortion -> SamplePortionUnit;
# N (new): Energy (joules); O (new): SamplePortion (mg)
$ws.Range("L1").Value = "SamplePortion"
$ws.Range("M1").Value = "SamplePortionUnit"
$ws.Range("N1").Value = "Energy"
$ws.Range("O1").Value = "SamplePortion"

# --- Row 2: type/unit annotations ---
$ws.Range("A2").Value = "#string"
$ws.Range("B2").Value = "#string"
$ws.Range("C2").Value = "#date"
$ws.Range("D2").Value = "#string"
$ws.Range("E2").Value = "#string"
$ws.Range("F2").Value = "#string"
$ws.Range("G2").Value = "#string"
$ws.Range("H2").Value = "#string"
$ws.Range("I2").Value = "#string"
$ws.Range("J2").Value = "#string"
$ws.Range("K2").Value = "#float"
$ws.Range("L2").Value = "#float,  unit:ml"
$ws.Range("M2").Value = "#string"
$ws.Range("N2").Value = "#float,  unit:joules"
$ws.Range("O2").Value = "#float,  unit:mg"

# --- Row 3: new French field-description row ---
$ws.Range("A3").Value = "#Manipulateur"
$ws.Range("B3").Value = "#Desc:IdentifiantEchantillon"
$ws.Range("C3").Value = "#Date"
$ws.Range("D3").Value = "#ModeOderatoireLaboratoire"
$ws.Range("E3").Value = "#AppareilLogicielCritique"
$ws.Range("F3").Value = "#ProduitCritique"
$ws.Range("G3").Value = "#LieuStockageDonneesBrutes"
$ws.Range("H3").Value = "#NumeroBombe"
$ws.Range("I3").Value = "#TypeCreuset"
$ws.Range("J3").Value = "#Programme"
$ws.Range("K3").Value = "#ValeurK"
$ws.Range("L3").Value = "#PriseEssai"
$ws.Range("M3").Value = "#UnitePriseEssai"
$ws.Range("N3").Value = "#Energie"
$ws.Range("O3").Value = "'"
$ws.Range("O3").ClearFormats()
